$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "DataFinal" to "var"
$ws.Name = "var"

# New PBI2016 (column H) values for rows 2..116, replacing the
# previous un-normalized (raw) GDP figures with values expressed
# in billions.
$hValues = @(19.36,11.86,160.0,101.1,557.5,10.55,1209.0,395.2,37.87,32.25,221.41,475.9,11.82,33.94,16.91,15.64,1795.7,53.78,10.87,2.96,20.16,32.62,1528.24,1.82,10.09,250.44,11233.28,282.82,51.6,195.09,75.7,99.94,332.93,24.19,24.0,240.61,2471.28,14.02,15.14,3466.79,55.01,195.22,3.5,7.97,21.72,127.51,2294.8,931.88,417.98,174.88,318.95,1875.8,14.07,4927.54,39.2,137.28,69.19,109.42,6.81,27.73,51.2,2.16,3.28,43.02,11.85,5.43,14.02,6.4,12.23,1077.9,11.19,4.38,103.31,11.94,10.66,21.18,783.53,188.22,13.29,10.29,368.82,65.48,278.65,57.91,20.76,36.05,191.9,318.63,472.04,206.29,188.49,1276.79,644.93,19.04,40.63,3.67,318.65,89.64,44.65,296.36,1232.08,82.4,51.77,515.65,671.31,49.77,37.13,18707.19,41.8,863.72,28.97,93.35,2694.28,20.95,20.55)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}
